# scan_example.xlsx edit script
# Adds a new "raw" field (mimeUri) to the model sheet, describing an
# extra scanned-output column, and updates the "properties" sheet's
# colOrder JSON array to include the new raw / raw_contentType /
# raw_uriFragment columns, matching how existing *_image0 columns are
# described in the model.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. "model" sheet: add header columns C1:H1 and a new data row A3:H3
#    describing the new "raw" mimeUri field (mirrors qrcode_image0 etc.)
# ---------------------------------------------------------------
$model = $wb.Worksheets.Item("model")

$model.Range("C1").Value = "elementType"
$model.Range("D1").Value = "properties.uriFragment.type"
$model.Range("E1").Value = "properties.uriFragment.elementType"
$model.Range("F1").Value = "properties.contentType.type"
$model.Range("G1").Value = "properties.contentType.elementType"
$model.Range("H1").Value = "properties.contentType.default"

$model.Range("A3").Value = "object"
$model.Range("B3").Value = "raw"
$model.Range("C3").Value = "mimeUri"
$model.Range("D3").Value = "string"
$model.Range("E3").Value = "rowpath"
$model.Range("F3").Value = "string"
$model.Range("G3").Value = "mimeType"
$model.Range("H3").Value = "application/json"

# Column widths for the newly-used columns (offset to compensate for the
# padding the engine adds between COM ColumnWidth and the stored XML width)
$model.Columns.Item(4).ColumnWidth = 24.998697916666668
$model.Columns.Item(5).ColumnWidth = 31.830729166666668
$model.Columns.Item(6).ColumnWidth = 33.166666666666664
$model.Columns.Item(7).ColumnWidth = 33.166666666666664
$model.Columns.Item(8).ColumnWidth = 28.330729166666668

$model.PageSetup.Orientation = 1

# ---------------------------------------------------------------
# 2. "properties" sheet: update colOrder array (row 2, column E) so it
#    includes the new raw / raw_contentType / raw_uriFragment entries
#    (alphabetically ordered, same as the rest of the list).
# ---------------------------------------------------------------
$props = $wb.Worksheets.Item("properties")

$colOrder = '["address","address_image0_contentType","address_image0_uriFragment","comments","comments_image0_contentType","comments_image0_uriFragment","fri_chores","fri_chores_image0_contentType","fri_chores_image0_uriFragment","mon_chores","mon_chores_image0_contentType","mon_chores_image0_uriFragment","name","name_image0_contentType","name_image0_uriFragment","qrcode","qrcode_image0_contentType","qrcode_image0_uriFragment","raw","raw_contentType","raw_uriFragment","roomNum","roomNum_image0_contentType","roomNum_image0_uriFragment","sat_chores","sat_chores_image0_contentType","sat_chores_image0_uriFragment","scan_output_directory","stay","stay_image0_contentType","stay_image0_uriFragment","sun_chores","sun_chores_image0_contentType","sun_chores_image0_uriFragment","thurs_chores","thurs_chores_image0_contentType","thurs_chores_image0_uriFragment","tues_chores","tues_chores_image0_contentType","tues_chores_image0_uriFragment","wed_chores","wed_chores_image0_contentType","wed_chores_image0_uriFragment"]'

$props.Range("E2").Value = $colOrder

# ---------------------------------------------------------------
# 3. Restore / update selections on each sheet to match the latest
#    interactive state, finishing on "properties" (the tab that was
#    active when the workbook was last saved).
# ---------------------------------------------------------------
$survey = $wb.Worksheets.Item("survey")
$survey.Range("B52").Select()

$choices = $wb.Worksheets.Item("choices")
$choices.Range("C15").Select()

$model.Range("A3").Select()

$settings = $wb.Worksheets.Item("settings")
$settings.Range("A3").Select()

$props.Range("E2").Select()
